# Update in read row, column and cell
#
# 1) The existing "Sheet1" gradebook gets one more row (A17 = "Hello World").
# 2) A brand new worksheet named "Sheet2" is inserted BEFORE "Sheet1" and
#    becomes the active tab; it holds a small 2x2 block (A1:B2) that all
#    reads "Hello World".

$wb = $excel.ActiveWorkbook

# --- Step 1: add the extra row to the existing gradebook sheet -------------
$gradeSheet = $wb.Worksheets.Item("Sheet1")
$gradeSheet.Range("A17").Value = "Hello World"
$gradeSheet.Range("A1").Select()

# --- Step 2: insert a new worksheet before Sheet1 and populate it ---------
$newSheet = $wb.Worksheets.Add($gradeSheet)
$newSheet.Name = "Sheet2"
$newSheet.Range("A1:B2").Value = "Hello World"
$newSheet.Range("A1").Select()
